$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000
$ws.Range("F2").Value = 3.15
$ws.Range("G2").Value = 4.7
$ws.Range("H2").Value = 1.97
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 1.94
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 1.75
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 1.75
$ws.Range("Q2").Value = 1.89
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.89
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.86
$ws.Range("V2").Value = 1.37
$ws.Range("W2").Value = 1.27
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("P3").Value = 1.87
$ws.Range("Q3").Value = 1.76
$ws.Range("AC4").Value = 8.6
$ws.Range("AE4").Value = 190
$ws.Range("AG4").Value = 10.5
$ws.Range("AJ4").Value = 22
$ws.Range("AK4").Value = 20
$ws.Range("F4").Value = 1.92
$ws.Range("G4").Value = 1.97
$ws.Range("J4").Value = 3.8
$ws.Range("N4").Value = 4.1
$ws.Range("R4").Value = 1.42
$ws.Range("T4").Value = 1.8
$ws.Range("U4").Value = 2.18
$ws.Range("Z4").Value = 34
$ws.Range("AB5").Value = 7.2
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 150
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 7.2
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 980
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 430
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 7
$ws.Range("H5").Value = 11
$ws.Range("J5").Value = 5.4
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 1.99
$ws.Range("Q5").Value = 1.96
$ws.Range("R5").Value = 1.35
$ws.Range("U5").Value = 1.63
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 980
$ws.Range("Z5").Value = 1000
$ws.Range("F6").Value = 2.58
$ws.Range("G6").Value = 3.7
$ws.Range("H6").Value = 2.26
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 4.8
$ws.Range("F7").Value = 1.59
$ws.Range("G7").Value = 1.95
$ws.Range("H7").Value = 2.04
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 980
$ws.Range("N7").Value = 1.94
$ws.Range("P7").Value = 1.94
$ws.Range("R7").Value = 1.08
$ws.Range("S7").Value = 1.63
$ws.Range("AB8").Value = 70
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 140
$ws.Range("AG8").Value = 40
$ws.Range("AI8").Value = 26
$ws.Range("AJ8").Value = 340
$ws.Range("AK8").Value = 140
$ws.Range("AO8").Value = 3.05
$ws.Range("F8").Value = 9.8
$ws.Range("H8").Value = 1.31
$ws.Range("I8").Value = 1.32
$ws.Range("K8").Value = 7.4
$ws.Range("Q8").Value = 1.3
$ws.Range("S8").Value = 1.74
$ws.Range("U8").Value = 2.56
$ws.Range("Y8").Value = 20
$ws.Range("F9").Value = 1.88
$ws.Range("G9").Value = 1.94
$ws.Range("J9").Value = 4.1
$ws.Range("K9").Value = 4.2
$ws.Range("P9").Value = 2.42
$ws.Range("R9").Value = 1.58
$ws.Range("U9").Value = 2.42
$ws.Range("S10").Value = 2.66
$ws.Range("U10").Value = 2.4
$ws.Range("AI11").Value = 510
$ws.Range("AN11").Value = 3.4
$ws.Range("H11").Value = 26
$ws.Range("T11").Value = 2.78
$ws.Range("U11").Value = 1.52
$ws.Range("X11").Value = 32
$ws.Range("G12").Value = 1.46
$ws.Range("H12").Value = 8.6
$ws.Range("AB13").Value = 27
$ws.Range("AC13").Value = 11.5
$ws.Range("AE13").Value = 18
$ws.Range("AH13").Value = 24
$ws.Range("AO13").Value = 7.8
$ws.Range("I13").Value = 1.56
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 2.3
$ws.Range("X13").Value = 980
$ws.Range("Y13").Value = 11
$ws.Range("Z13").Value = 12
$ws.Range("G14").Value = 1.54
$ws.Range("H14").Value = 1.04
$ws.Range("I14").Value = 15
$ws.Range("K14").Value = 980
$ws.Range("N14").Value = 1.26
$ws.Range("P14").Value = 1.25
$ws.Range("Q14").Value = 1.18
$ws.Range("R14").Value = 1.18
$ws.Range("S14").Value = 1.18
$ws.Range("V14").Value = 1.07
